$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 3

# Row 5
$ws.Range("B5").Value = 6

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 8

# Row 7
$ws.Range("B7").Value = 8

# Row 8
$ws.Range("B8").Value = 5

# Row 9
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = 6

# Row 10
$ws.Range("A10").Value = 3
$ws.Range("B10").Value = 5

# Row 11
$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 8
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1.375
